$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 and 27 swap: Monero/Stellar order flipped with updated price/volume data
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1277"
$ws.Range("E26").Value = "  +4.12%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "140.44"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.821.47"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.84"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.56"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5076"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06402"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.35"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07807"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.260"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.24"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.855.06"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5588"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.30"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7516"
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.837.62"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.17"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.299"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.809"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.998"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.835"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.738"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.40"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.239"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04866"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.288"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.183"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.557"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.379"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8959"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.132.30"
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.547"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5465"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01559"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.571"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7955"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.20"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.779.92"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -6.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4441"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.98"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05059"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.572"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.37%  "
